$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "DKS / Desktop Computer" device-type rows (rows 8-10) -- the
# three rows below the header shift up, the dimension shrinks to G16, and
# the now-unused shared strings for that block drop out of the table.
$ws.Range("A8:A10").EntireRow.Delete()

# Leave the active cell where it landed after the edit.
$ws.Range("E10").Select()

# Page setup touched during the same save (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
